$wb = $excel.ActiveWorkbook

# Mapping of row -> new "want to go" count (column F) for sheets 展览 and 全部类型
$updates = @{
    2  = 8007
    3  = 7608
    10 = 150
    11 = 223
    12 = 688
    13 = 110
    14 = 1203
    15 = 57
    16 = 43
    17 = 7
    19 = 99
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
